$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Đơn phụ phẫu 2")

# Append a new blank data row (row 2) below the header row, mirroring the
# report template's per-column types: numeric "total" columns get 0,
# text columns get an (empty) value. Touching a format property on the
# text cells forces them to be materialised in the sheet (so the row/
# dimension grows to A1:T2) even though their value is blank.
$ws.Range("A2").Font.Bold = $false
$ws.Range("B2").Value = 0
$ws.Range("C2").Font.Bold = $false
$ws.Range("D2").Font.Bold = $false
$ws.Range("E2").Font.Bold = $false
$ws.Range("F2").Font.Bold = $false
$ws.Range("G2").Font.Bold = $false
$ws.Range("H2").Font.Bold = $false
$ws.Range("I2").Value = 0
$ws.Range("J2").Font.Bold = $false
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Font.Bold = $false
$ws.Range("R2").Font.Bold = $false
$ws.Range("S2").Font.Bold = $false
$ws.Range("T2").Font.Bold = $false
